$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: row index, A-value (as string token so "asd15" stays text), B-value (URL or empty)
$rows = @(
    [PSCustomObject]@{ Row = 2; A = 1; ANumeric = $true; B = 'https://masterwatt.ru/catalog/komplektuyushchie-i-zapchasti-k-kotlam/komnatnyy-regulyator-thermolink-b-ebus-24-v-montazh-naruzhny/' },
    [PSCustomObject]@{ Row = 3; A = 2; ANumeric = $true; B = 'https://masterwatt.ru/catalog/komplektuyushchie-i-zapchasti-k-kotlam/komnatnyy-regulyator-thermolink-b-ebus-24-v-montazh-naruzhny/' },
    [PSCustomObject]@{ Row = 4; A = 3; ANumeric = $true; B = 'https://masterwatt.ru/catalog/komplektuyushchie-i-zapchasti-k-kotlam/komnatnyy-regulyator-thermolink-b-ebus-24-v-montazh-naruzhny/' },
    [PSCustomObject]@{ Row = 5; A = 4; ANumeric = $true; B = 'https://masterwatt.ru/catalog/komplektuyushchie-i-zapchasti-k-kotlam/komnatnyy-regulyator-thermolink-b-ebus-24-v-montazh-naruzhny/' },
    [PSCustomObject]@{ Row = 6; A = 5; ANumeric = $true; B = 'https://masterwatt.ru/catalog/komplektuyushchie-i-zapchasti-k-kotlam/komnatnyy-regulyator-thermolink-b-ebus-24-v-montazh-naruzhny/' },
    [PSCustomObject]@{ Row = 7; A = 6; ANumeric = $true; B = 'https://masterwatt.ru/catalog/komplektuyushchie-i-zapchasti-k-kotlam/komnatnyy-regulyator-thermolink-b-ebus-24-v-montazh-naruzhny/' },
    [PSCustomObject]@{ Row = 8; A = 7; ANumeric = $true; B = 'https://masterwatt.ru/catalog/komplektuyushchie-i-zapchasti-k-kotlam/komplekt-perevoda-na-szhizhennyy-gaz-dlya-kotla-gepard-12-i-/' },
    [PSCustomObject]@{ Row = 9; A = 8; ANumeric = $true; B = 'https://masterwatt.ru/catalog/prinadlezhnosti-dlya-sistem-podderzhaniya-davleniya-i-podpit/plata-plavnogo-puska-dvukh-nasosov-sanftanlaufplatine-2-pump/' },
    [PSCustomObject]@{ Row = 10; A = 9; ANumeric = $true; B = 'https://masterwatt.ru/catalog/komplektuyushchie-i-zapchasti-k-kotlam/komnatnyy-regulyator-thermolink-b-ebus-24-v-montazh-naruzhny/' },
    [PSCustomObject]@{ Row = 11; A = 10; ANumeric = $true; B = 'https://masterwatt.ru/catalog/prinadlezhnosti-dlya-sistem-podderzhaniya-davleniya-i-podpit/datchik-urovnya-zhidkosti-v-emkostyakh-vg-rg-gg-4-20-ma-0-10~1/' },
    [PSCustomObject]@{ Row = 12; A = 11; ANumeric = $true; B = 'https://masterwatt.ru/catalog/komplektuyushchie-i-zapchasti-k-kotlam/komplekt-perevoda-na-szhizhennyy-gaz-dlya-kotla-gepard-12-i-/' },
    [PSCustomObject]@{ Row = 13; A = 12; ANumeric = $true; B = 'https://masterwatt.ru/catalog/komplektuyushchie-i-zapchasti-k-kotlam/komplekt-perevoda-na-szhizhennyy-gaz-dlya-kotla-pantera-12-k/' },
    [PSCustomObject]@{ Row = 14; A = 13; ANumeric = $true; B = 'https://masterwatt.ru/catalog/vodorozetki/ugolnik-ustanovochnyy-gx-20x1-2-f-45-mm/' },
    [PSCustomObject]@{ Row = 15; A = 14; ANumeric = $true; B = 'https://masterwatt.ru/catalog/rasshiritelnye-baki/rasshiritelnyy-bak-de-6000-s-16-dn-65-pn-16-siniy-c-t-v-cert/' },
    [PSCustomObject]@{ Row = 16; A = 'asd15'; ANumeric = $false; B = '' },
    [PSCustomObject]@{ Row = 17; A = 16; ANumeric = $true; B = 'https://masterwatt.ru/catalog/komplektuyushchie-i-zapchasti-k-kotlam/komplekt-perevoda-na-szhizhennyy-gaz-dlya-kotla-gepard-12-i-/' },
    [PSCustomObject]@{ Row = 18; A = 17; ANumeric = $true; B = 'https://masterwatt.ru/catalog/komplektuyushchie-i-zapchasti-k-kotlam/datchik-boylera-ntc-dlya-kotlov-pantera-12-kto-25-kto-25-koo/' },
    [PSCustomObject]@{ Row = 19; A = 18; ANumeric = $true; B = 'https://masterwatt.ru/catalog/separatory/separator-gryazi-i-shlama-reflex-exdirt-r-hc-s-d-200-16-bar-/' },
    [PSCustomObject]@{ Row = 20; A = 19; ANumeric = $true; B = 'https://masterwatt.ru/catalog/komplektuyushchie-i-zapchasti-k-kotlam/komnatnyy-regulyator-thermolink-b-ebus-24-v-montazh-naruzhny/' },
    [PSCustomObject]@{ Row = 21; A = 20; ANumeric = $true; B = 'https://masterwatt.ru/catalog/komplektuyushchie-i-zapchasti-k-kotlam/komnatnyy-regulyator-thermolink-b-ebus-24-v-montazh-naruzhny/' },
    [PSCustomObject]@{ Row = 22; A = 21; ANumeric = $true; B = 'https://masterwatt.ru/catalog/komplektuyushchie-i-zapchasti-k-kotlam/komplekt-perevoda-na-szhizhennyy-gaz-dlya-kotla-gepard-12-i-/' },
    [PSCustomObject]@{ Row = 23; A = 22; ANumeric = $true; B = 'https://masterwatt.ru/catalog/prinadlezhnosti-dlya-sistem-podderzhaniya-davleniya-i-podpit/vakuumnyy-deaerator-servitec-60-t-control-touch/' },
    [PSCustomObject]@{ Row = 24; A = 23; ANumeric = $true; B = 'https://masterwatt.ru/catalog/komplektuyushchie-i-zapchasti-k-kotlam/komplekt-perevoda-na-szhizhennyy-gaz-dlya-kotla-gepard-12-i-/' },
    [PSCustomObject]@{ Row = 25; A = 1; ANumeric = $true; B = 'https://masterwatt.ru/catalog/komplektuyushchie-i-zapchasti-k-kotlam/komnatnyy-regulyator-thermolink-b-ebus-24-v-montazh-naruzhny/' },
    [PSCustomObject]@{ Row = 26; A = 2; ANumeric = $true; B = 'https://masterwatt.ru/catalog/komplektuyushchie-i-zapchasti-k-kotlam/komnatnyy-regulyator-thermolink-b-ebus-24-v-montazh-naruzhny/' },
    [PSCustomObject]@{ Row = 27; A = 3; ANumeric = $true; B = 'https://masterwatt.ru/catalog/komplektuyushchie-i-zapchasti-k-kotlam/komnatnyy-regulyator-thermolink-b-ebus-24-v-montazh-naruzhny/' },
    [PSCustomObject]@{ Row = 28; A = 4; ANumeric = $true; B = 'https://masterwatt.ru/catalog/komplektuyushchie-i-zapchasti-k-kotlam/komnatnyy-regulyator-thermolink-b-ebus-24-v-montazh-naruzhny/' },
    [PSCustomObject]@{ Row = 29; A = 5; ANumeric = $true; B = 'https://masterwatt.ru/catalog/komplektuyushchie-i-zapchasti-k-kotlam/komnatnyy-regulyator-thermolink-b-ebus-24-v-montazh-naruzhny/' },
    [PSCustomObject]@{ Row = 30; A = 6; ANumeric = $true; B = 'https://masterwatt.ru/catalog/komplektuyushchie-i-zapchasti-k-kotlam/komnatnyy-regulyator-thermolink-b-ebus-24-v-montazh-naruzhny/' },
    [PSCustomObject]@{ Row = 31; A = 7; ANumeric = $true; B = 'https://masterwatt.ru/catalog/komplektuyushchie-i-zapchasti-k-kotlam/komplekt-perevoda-na-szhizhennyy-gaz-dlya-kotla-gepard-12-i-/' },
    [PSCustomObject]@{ Row = 32; A = 8; ANumeric = $true; B = 'https://masterwatt.ru/catalog/prinadlezhnosti-dlya-sistem-podderzhaniya-davleniya-i-podpit/plata-plavnogo-puska-dvukh-nasosov-sanftanlaufplatine-2-pump/' }
)

foreach ($r in $rows) {
    if ($r.ANumeric) {
        $ws.Cells.Item($r.Row, 1).Value = [double]$r.A
    } else {
        $ws.Cells.Item($r.Row, 1).Value = $r.A
    }
    $ws.Cells.Item($r.Row, 2).Value = $r.B
}

Write-Output "done"